# "Capitulo I, parte Ashley"
# Slide 1, subtitle placeholder ("Subtítulo 2"): split the run so the
# space between "Yamilka" and "Gómez, ..." becomes its own run at a
# smaller size (1600), and "Yamilka" no longer carries an explicit
# run-level size override (it already visually matches the layout's
# default size).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# 1) Drop "Yamilka"'s explicit sz="2400" override: delete the word and
#    retype it so the new run inherits formatting instead of keeping the
#    old explicit size.
$tr.Characters(1, 7).Text = ""
$tr.InsertBefore("Yamilka") | Out-Null

# 2) Give the single space that follows "Yamilka" its own run sized
#    down to 1600 (splits it off from the "Gómez, ..." run).
$tr.Characters(8, 1).Font.Size = 16
